$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B5").Value = "Cuối kỳ"
$ws.Range("B5").Select()
